# RV25: XGBoost + variables Intento de Red Neuronal
# Adds a new results row (row 23 / "No." 25) to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recipeText = "update_role(property_type,  area, dist_nearest_restaurant,`n                dist_nearest_parques, baños, n_pisos_numerico,dist_nearest_universidades,`n                terraza, ascensor, estrato, dis_centro, dis_andino,`n                restaurant_1km, parques_1km, discotecas_1km, colegios_1km,`n                iluminado, exterior, remodelado, restaurant_300m, parques_300m,`n                discotecas_300m, colegios_300m, restaurant_100m, parques_100m, `n                discotecas_100m, colegios_100m,"

$covarText = "property_type,  area, dist_nearest_restaurant,`n                dist_nearest_parques, baños, n_pisos_numerico,dist_nearest_universidades,`n                terraza, ascensor, estrato, dis_centro, dis_andino,`n                restaurant_1km, parques_1km, discotecas_1km, colegios_1km,`n                iluminado, exterior, remodelado, restaurant_300m, parques_300m,`n                discotecas_300m, colegios_300m, restaurant_100m, parques_100m, `n                discotecas_100m, colegios_100m,"

$row = 23

$ws.Cells.Item($row, 1).Value = 25
$ws.Cells.Item($row, 2).Value = $recipeText
$ws.Cells.Item($row, 3).Value = $covarText
$ws.Cells.Item($row, 4).Value = "chapitrain"
$ws.Cells.Item($row, 5).Value = "XGBoost"
$ws.Cells.Item($row, 6).Value = "trees: 5000; mtry: 59; tree depth: 4"
$ws.Cells.Item($row, 7).Value = 63366019
$ws.Cells.Item($row, 8).Value = 199574167

# Match the "Comma" cell style already used by the G/H columns.
$ws.Range($ws.Cells.Item($row, 7), $ws.Cells.Item($row, 8)).NumberFormat = '_-* #,##0_-;\-* #,##0_-;_-* "-"??_-;_-@_-'

# The wrapped recipe/covariable cells pick up a left-aligned style.
$ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 3)).HorizontalAlignment = 1

# Typing a multi-line value bumps the row to a custom height; auto-fit puts
# it back in line with every other (also multi-line) row in the sheet.
$ws.Rows.Item($row).AutoFit()

# Leave the selection on the last entry typed, as the author did.
[void]$ws.Range("H23").Select()

Write-Output "Row $row appended"
